# Update the "想去人数" (F) and "最低票价" (G) figures for the latest data pull.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1465
$ws1.Range("F4").Value = 1754
$ws1.Range("F5").Value = 32
$ws1.Range("F6").Value = 143
$ws1.Range("F10").Value = 552
$ws1.Range("F11").Value = 25
$ws1.Range("F12").Value = 78
$ws1.Range("G12").Value = 50
$ws1.Range("F13").Value = 146
$ws1.Range("F14").Value = 22
$ws1.Range("F16").Value = 69
$ws1.Range("F17").Value = 103
$ws1.Range("F18").Value = 4666
$ws1.Range("F19").Value = 43
$ws1.Range("F20").Value = 819
$ws1.Range("F21").Value = 101
$ws1.Range("F22").Value = 2194
$ws1.Range("F24").Value = 16
$ws1.Range("F25").Value = 2057

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 75

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1465
$ws4.Range("F4").Value = 1754
$ws4.Range("F5").Value = 32
$ws4.Range("F6").Value = 143
$ws4.Range("F10").Value = 552
$ws4.Range("F11").Value = 25
$ws4.Range("F12").Value = 78
$ws4.Range("G12").Value = 50
$ws4.Range("F13").Value = 146
$ws4.Range("F14").Value = 22
$ws4.Range("F16").Value = 69
$ws4.Range("F17").Value = 103
$ws4.Range("F18").Value = 4666
$ws4.Range("F19").Value = 75
$ws4.Range("F20").Value = 43
$ws4.Range("F22").Value = 819
$ws4.Range("F23").Value = 101
$ws4.Range("F24").Value = 2194
$ws4.Range("F26").Value = 16
$ws4.Range("F27").Value = 2057
